# Changing from job-aids to Deliverables
#
# 1) The "Can 16" shape inside the "Group 3" group on slide 1 has its
#    label text changed from "Job aids" to "Deliverables".
# 2) The containing group ("Group 3") is nudged up slightly (its Top
#    changes while the child-coordinate space (chOff) stays put).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Locate the top-level group shape on the slide -----------------
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Group 3") {
        $grp = $cand
    }
}
if ($grp -eq $null) {
    $grp = $s.Shapes.Item(1)
}

# --- Reposition the group: y offset 488651 EMU -> 476672 EMU -------
# PowerPoint COM works in points; 1 pt = 12700 EMU.
$emuPerPoint = 12700
$grp.Top = 476672 / $emuPerPoint

# --- Rename the "Job aids" can to "Deliverables" --------------------
$target = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($item.HasTextFrame -and $item.TextFrame.TextRange.Text -eq "Job aids") {
        $target = $item
    }
}
if ($target -eq $null) {
    $target = $grp.GroupItems.Item("Can 16")
}

$target.TextFrame.TextRange.Text = "Deliverables"
